# Cronograma.xlsx edit script
# Adds a "Feito" (Done) column to the Plan1 schedule sheet, marks several
# tasks as done (green "x"), and makes Plan1 the active/selected sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Plan1")

# Green color used for the "done" marks (RGB 0,176,80 == 00B050)
$doneGreen = 5287936

# Slide the old column D (the C2/B2 percentage formula) over to column E,
# and give the freed-up column D the plain look of column C, matching
# what a real column insert before D would produce.
$ws1.Range("E2").NumberFormat = $ws1.Range("D2").NumberFormat
$ws1.Range("E2").Formula = $ws1.Range("D2").Formula
$ws1.Range("D2").Formula = $null
$ws1.Range("C2").Copy()
$ws1.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Column D header and width.
$ws1.Range("D1").Value = "Feito"
$ws1.Columns.Item(4).ColumnWidth = 4.25

# Rows whose tasks are marked as finished ("x" in column D, green text
# across the row). Row 9 is marked as "done" styling-wise but with no
# visible mark.
$doneRowsWithMark = @(5, 7, 8, 12, 16)
foreach ($r in $doneRowsWithMark) {
    $ws1.Range("D$r").Value = "x"
}
$ws1.Range("D9").Value = ""

$doneRows = @(5, 7, 8, 9, 12, 16)
foreach ($r in $doneRows) {
    $rowRange = $ws1.Range("A" + $r + ":D" + $r)
    $rowRange.Font.Color = $doneGreen
}

# Make Plan1 the active sheet/tab, with B22 selected (matches the
# workbook's last saved view).
$ws1.Activate()
$ws1.Range("B22").Select()
